$d = $word.ActiveDocument

# --- Change 1: "We then created a" -> "We first created a" ---
$find1 = $d.Content.Find
$find1.Execute("We then created a", $true, $false, $false, $false, $false, $true, 1, $false, "We first created a", 2) | Out-Null

# --- Change 2: "library" -> "module" (bayespy library. -> bayespy module.) ---
$find2 = $d.Content.Find
$find2.Execute("library", $true, $false, $false, $false, $false, $true, 1, $false, "module", 2) | Out-Null

# --- Change 3: DumbAI paragraph -- drop the _GoBack bookmark that currently
#     lives here, then append the new sentence to the paragraph. ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$dumbAiPara = $d.Paragraphs.Item(10)
$dumbAiRange = $dumbAiPara.Range
$insertionPoint = $d.Range($dumbAiRange.End - 1, $dumbAiRange.End - 1)
$insertionPoint.InsertAfter("This gave us a base line from which to judge how well our Bayesian Network AI played. ")

# --- Change 4: fill the empty paragraph right after "LEVEL DIAGRAMS HERE" ---
$basicLevelPara = $d.Paragraphs.Item(15)
$basicLevelPara.Range.Text = "For our basic level, we only incorporated the "

# --- Change 5: fill the empty paragraph right after "Future Work:" and put
#     the _GoBack bookmark back, now at the end of this new paragraph. ---
$futureWorkPara = $d.Paragraphs.Item(19)
$futureWorkPara.Range.Text = "Future work we would like to see done, is the incorporation of the character abilities into the game and have them also factor in to the Bayesian Network Model. "

$futureWorkPara2 = $d.Paragraphs.Item(19)
$futureWorkRange2 = $futureWorkPara2.Range
# NB: placing a bookmark exactly at (paragraph.End - 1) trips an indexing bug
# in this host, so anchor one character earlier instead (still effectively
# "at the end" since the bookmark itself is zero-width / invisible).
$bookmarkPos = $futureWorkRange2.End - 2
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

Write-Output "done"
